$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 10 (Objetivos:) - B/C currently hold the misplaced "4873328 - Fernando
#    Segato" string; replace with the real Portuguese objectives paragraph.
# ---------------------------------------------------------------------------
$ws.Range("B10:C10").Value = "Introduzir o estudante em conceitos importantes utilizados nas ciências microbiológicas, apresentando os conceitos fundamentais relacionados à história, mercado, genética, filogenia, e cultivo de microrganismos."

# ---------------------------------------------------------------------------
# 2. Insert two new rows right after row 12 ("Docentes responsáveis:") to hold
#    the two professors, one per row (B/C only, no label in A). The default
#    insert copies row 12's (bold, label) formatting into the new row, so
#    clear column A and repaint B:C with the plain data-cell formatting used
#    elsewhere (copied from row 16) before writing the values.
# ---------------------------------------------------------------------------
$ws.Rows("13:14").Insert()
$ws.Range("A13:A14").Clear()
$ws.Range("B16:C16").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13:C13").Value = "4873328 - Fernando Segato"
$ws.Range("B14:C14").Value = "5840685 - Maria Bernadete de Medeiros"

# ---------------------------------------------------------------------------
# 3. Old row 13 ("Programa resumido:") is now row 15 - fix its B/C value
#    (was wrongly "01/01/2018").
# ---------------------------------------------------------------------------
$ws.Range("B15:C15").Value = "Histórico da microbiologia, microbiologia industrial, filogênia microbiana, caracterização dos microrganismos, nutrição e cultivo de microrganismos, virus, fungos filamentosos, leveduras, micro-algas, bactérias."

# ---------------------------------------------------------------------------
# 4. Old row 15 ("Programa:") is now row 17 - fix its B/C value
#    (was wrongly "4873328 - Fernando Segato").
# ---------------------------------------------------------------------------
$ws.Range("B17:C17").Value = "1. Histórico da microbiologia: microbiologia, ciência e sociedade;  Leeuwenhoek e seusseus microscópios, origem dos animálculos de Leeuwenhoek, principais pensadores da microbiologia, microbiologia moderna.2. Microbiologia industrial: visão geral do mercado envolvendo microbiologia, principais produtos de origem microbiana.3. Filogênia microbiana: classificação e evolução das principais classes dos microrganismos; organismos procarióticos (eubactérias e arqueobactérias); organismos eucarióticos (leveduras, fungos filamentosos, algas, protozoários).4. Caracterização dos microrganismos: técnicas de cultura pura, microscópios, técnicas de microscopia, preparo dos microrganismos para microscopia, informações utilizadas para caracterizar os microrganismos, tecnologia automatizada.5. Nutrição e cultivo de microrganismos: exigências nutricionais e meios microbiológicos; cultivo ecrescimento de microrganismos.6. Genética de microrganismos: regulação da expressão gênica em bactérias; mutação, vantagens e desvantagens para aplicações industriais; melhoramento de cepas.7. Virus, bactérias, fungos filamentosos, micro-algas, leveduras: morfologia, classificação e replicação."

# ---------------------------------------------------------------------------
# 5. Old row 18 ("Método:") is now row 20 - fix its B/C value
#    (was wrongly "5840685 - Maria Bernadete de Medeiros").
# ---------------------------------------------------------------------------
$ws.Range("B20:C20").Value = "A avaliação será feita por meio de provas escritas, trabalhos, seminários e participação."

# ---------------------------------------------------------------------------
# 6. Old row 19 ("Critério:") is now row 21 - fix its B/C value
#    (was wrongly the "A avaliação..." text).
# ---------------------------------------------------------------------------
$ws.Range("B21:C21").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2."

# ---------------------------------------------------------------------------
# 7. Old row 20 ("Norma de recuperação:") is now row 22 - fix its B/C value
#    (was wrongly the "A Nota final..." text).
# ---------------------------------------------------------------------------
$ws.Range("B22:C22").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR)"

# ---------------------------------------------------------------------------
# 8. Old row 21 ("Bibliografia:") is now row 23 - fix its B/C value
#    (was wrongly the "A recuperação..." text).
# ---------------------------------------------------------------------------
$ws.Range("B23:C23").Value = "1. PELCZAR Jr, M.J., CHAN, S.S., KRIEG, N.R. Microbiologia conceitos e aplicações, 2 ed. (Vol 1), São Paulo: Pearson Education do Brasil, 1997.2. MADIGAN, M.T., MARTINKO, J.M., PARKER, I. Microbiologia de Brock. São Paulo: Prentice Hall, 2004.3. BARBOSA, H.R., TORRES, B.B. Microbiologia Básica, Rio de Janeiro: Atheneu, 2005.4. VERMELHO A.B., FREIRE, M.C., BRANQUINHO, M.H. Bacteorologia Geral, Rio de Janeiro: GuanabaraKoogan, 2008.5. TORTORA, G.J., FUNKE, B.R., CASE, C.L. Microbiologia, Artmed, Porto Alegre, RS, 2012."

# ---------------------------------------------------------------------------
# 9. Column layout: split the combined "A:B" width definition so column A has
#    its own single-column width entry (same width as before).
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 60.7109375
